# Scheduled-runner style refresh of cached Leve profit figures across the
# crafting-job sheets (currentAveragePrice / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns H, I, J, K, L, M, N).
# Some rows' NQ profit (column M) is cleared entirely rather than zeroed,
# matching the source data's convention of omitting a cell when there is
# no computed profit for that row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2704
$ws.Range("I15").Value = 2704
$ws.Range("K15").Value = 8112
$ws.Range("M15").Value = -7943

$ws.Range("H17").Value = 2815.2
$ws.Range("I17").Value = 933.3333
$ws.Range("K17").Value = 2799.9999
$ws.Range("M17").Value = -2631.9999

$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 1333.3334
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 4000.0002
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -4540.0002

$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 1333.3334
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 4000.0002
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -5872.0002

$ws.Range("H118").Value = 467.4
$ws.Range("I118").Value = 467.4
$ws.Range("K118").Value = 1402.2
$ws.Range("M118").Value = 254.8000000000002

$ws.Range("H132").Value = 2084.3333
$ws.Range("I132").Value = 2084.3333
$ws.Range("K132").Value = 6252.999899999999
$ws.Range("M132").Value = -3722.999899999999

$ws.Range("H137").Value = 2733.3333
$ws.Range("I137").Value = 2630.4736
$ws.Range("J137").Value = 2911
$ws.Range("K137").Value = 7891.4208
$ws.Range("L137").Value = 8733
$ws.Range("M137").Value = -5341.4208
$ws.Range("N137").Value = -13833

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 59800
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = $null

$ws.Range("H32").Value = 17711.736
$ws.Range("I32").Value = 17711.736
$ws.Range("K32").Value = 17711.736
$ws.Range("M32").Value = -17424.736

$ws.Range("H61").Value = 1900
$ws.Range("I61").Value = 1900
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1900
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = -1688

$ws.Range("H74").Value = 12730.488
$ws.Range("I74").Value = 12997.639
$ws.Range("K74").Value = 12997.639
$ws.Range("M74").Value = -12123.639

$ws.Range("H77").Value = 12730.488
$ws.Range("I77").Value = 12997.639
$ws.Range("K77").Value = 64988.19499999999
$ws.Range("M77").Value = -60620.19499999999

$ws.Range("H99").Value = 59800
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null

$ws.Range("H102").Value = 3401.25
$ws.Range("I102").Value = 3401.25
$ws.Range("K102").Value = 3401.25
$ws.Range("M102").Value = -1779.25

$ws.Range("H110").Value = 5536.25
$ws.Range("I110").Value = 2755.7144
$ws.Range("K110").Value = 2755.7144
$ws.Range("M110").Value = -710.7143999999998

$ws.Range("H122").Value = 3569.6
$ws.Range("I122").Value = 3569.6
$ws.Range("K122").Value = 10708.8
$ws.Range("M122").Value = -8258.799999999999

$ws.Range("H136").Value = 1900
$ws.Range("I136").Value = 1900
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5700
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -3150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8147.25
$ws.Range("I86").Value = 3989
$ws.Range("J86").Value = 9533.333000000001
$ws.Range("K86").Value = 3989
$ws.Range("L86").Value = 9533.333000000001
$ws.Range("M86").Value = -2866
$ws.Range("N86").Value = -11779.333

$ws.Range("H89").Value = 8147.25
$ws.Range("I89").Value = 3989
$ws.Range("J89").Value = 9533.333000000001
$ws.Range("K89").Value = 19945
$ws.Range("L89").Value = 47666.665
$ws.Range("M89").Value = -14329
$ws.Range("N89").Value = -58898.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1922.3182
$ws.Range("I31").Value = 1806.3334
$ws.Range("K31").Value = 1806.3334
$ws.Range("M31").Value = -1511.3334

$ws.Range("H34").Value = 1922.3182
$ws.Range("I34").Value = 1806.3334
$ws.Range("K34").Value = 1806.3334
$ws.Range("M34").Value = -1604.3334

$ws.Range("H58").Value = 4171.696
$ws.Range("I58").Value = 3759.476
$ws.Range("K58").Value = 3759.476
$ws.Range("M58").Value = -3556.476

$ws.Range("H132").Value = 4249.75
$ws.Range("I132").Value = 4280.2
$ws.Range("J132").Value = 4199
$ws.Range("K132").Value = 12840.6
$ws.Range("L132").Value = 12597
$ws.Range("M132").Value = -10310.6
$ws.Range("N132").Value = -17657

$ws.Range("H136").Value = 4171.696
$ws.Range("I136").Value = 3759.476
$ws.Range("K136").Value = 11278.428
$ws.Range("M136").Value = -8728.428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5412
$ws.Range("I113").Value = 5412
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5412
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -3242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7789.615
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590

$ws.Range("H27").Value = 7789.615
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214

$ws.Range("H93").Value = 2306
$ws.Range("I93").Value = 1834
$ws.Range("J93").Value = 3250
$ws.Range("K93").Value = 1834
$ws.Range("L93").Value = 3250
$ws.Range("M93").Value = -586
$ws.Range("N93").Value = -5746

$ws.Range("H136").Value = 2940.9285
$ws.Range("I136").Value = 2874.7273
$ws.Range("J136").Value = 3183.6667
$ws.Range("K136").Value = 8624.1819
$ws.Range("L136").Value = 9551.000100000001
$ws.Range("M136").Value = -6074.1819
$ws.Range("N136").Value = -14651.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4750.4287
$ws.Range("I122").Value = 4750.4287
$ws.Range("K122").Value = 14251.2861
$ws.Range("M122").Value = -11801.2861

$ws.Range("H136").Value = 2082.5417
$ws.Range("I136").Value = 2316.2104
$ws.Range("K136").Value = 6948.6312
$ws.Range("M136").Value = -4398.6312
